$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "SP23092022232857"
$ws.Range("B21").Value = 280
$ws.Range("C21").Value = "SIGMACORT CREAM 1% 50G"
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 799
